# "semana 28 de 2025" - add week 28 column (AE) to the weekly IRA hospital report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new week-28 column. Use a leading apostrophe so Excel stores
# the numeric-looking label "28" as text, matching the other week headers
# (1, 2, 3, ... 27) already on row 1.
$ws.Range("AE1").Value = "'28"

# Week-28 counts for each UPGD row. Only rows that already carry data for
# previous weeks receive a value here (the remaining rows stay blank, same
# as the gaps already present in other week columns).
$aeValues = @{
  2  = 0
  4  = 0
  5  = 0
  6  = 23
  7  = 0
  8  = 15
  10 = 0
  12 = 0
  13 = 0
  14 = 0
  15 = 0
  16 = 0
  17 = 0
  18 = 0
  22 = 0
  23 = 0
  24 = 0
  25 = 4
  26 = 0
  28 = 163
  29 = 0
  30 = 10
  31 = 0
  32 = 0
  34 = 1
  35 = 18
  36 = 0
  37 = 0
  38 = 0
  39 = 0
  40 = 0
  41 = 0
  42 = 0
  43 = 0
  44 = 0
  45 = 0
  46 = 0
  47 = 0
  48 = 0
  49 = 0
  50 = 0
  52 = 0
  53 = 0
  54 = 0
  55 = 0
  56 = 0
  57 = 0
}

foreach ($row in $aeValues.Keys) {
  $ws.Cells.Item($row, 31).Value = $aeValues[$row]
}

# Row 28 (CLINICA LOS ROSALES) was also missing its week-26 (AC) value;
# fill it in alongside the new week-28 figure.
$ws.Range("AC28").Value = 0
